$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capacità di trasmissione MW")

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
